{"js": "// The title paragraph originally reads:\n//   \"RENCANA PROGRAM AUDIT INTERNAL PUSKESMAS\"\n// It is revised to:\n//   \"RENCANA AUDIT INTERNAL PUSKESMAS\"\n// (the word \"PROGRAM \" is removed), and the run is split in two at the\n// edit point, with Word's \"_GoBack\" last-edit-position bookmark landing\n// there. The \"_GoBack\" bookmark used to sit further down the document\n// (right before the \"Metoda audit:\" run); that stale copy is removed.\n\nconst body = context.document.body;\n\n// 1) Drop the old \"_GoBack\" bookmark (wherever Word last left it).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Remove the word \"PROGRAM \" from the title.\nconst progRanges = body.search(\"PROGRAM \", { matchCase: true });\nprogRanges.load(\"items\");\nawait context.sync();\n\nif (progRanges.items.length > 0) {\n  progRanges.items[0].delete();\n  await context.sync();\n}\n\n// 3) Re-plant \"_GoBack\" right after \"RENCANA \" (i.e. between the two\n//    words that now make up the title), splitting the run in the process.\nconst renRanges = body.search(\"RENCANA \", { matchCase: true });\nrenRanges.load(\"items\");\nawait context.sync();\n\nif (renRanges.items.length > 0) {\n  const afterRencana = renRanges.items[0].getRange(\"End\");\n  afterRencana.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The title paragraph originally reads:\n#   \"RENCANA PROGRAM AUDIT INTERNAL PUSKESMAS\"\n# It is revised to:\n#   \"RENCANA AUDIT INTERNAL PUSKESMAS\"\n# (the word \"PROGRAM \" is removed). Word's \"_GoBack\" last-edit-position\n# bookmark, which used to sit further down the document (right before the\n# \"Metoda audit:\" run), is moved to land at the edit point in the title -\n# i.e. right after \"RENCANA \".\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark (wherever Word last left it).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Remove the word \"PROGRAM \" from the title.\n$progRange = $d.Range()\n$null = $progRange.Find.Execute(\"PROGRAM \")\n$progRange.Text = \"\"\n\n# 3) Re-plant \"_GoBack\" right after \"RENCANA \" (i.e. between the two\n#    words that now make up the title).\n$renRange = $d.Range()\n$null = $renRange.Find.Execute(\"RENCANA \")\n$bmRange = $d.Range($renRange.End, $renRange.End)\n$null = $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
